# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.920.41"
$ws.Range("E2").Value = "  +10.56%  "
$ws.Range("D3").Value = "1.811.22"
$ws.Range("E3").Value = "  +7.57%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'227.73"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").Value = "'0.541"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'30.91"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "'47.83"
$ws.Range("E9").Value = "  +7.82%  "
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("D11").Value = "'0.0665"
$ws.Range("E11").Value = "  +6.43%  "
$ws.Range("D12").Value = "'0.0929"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "2.073.08"
$ws.Range("E13").Value = "  +7.52%  "
$ws.Range("D14").Value = "1.814.47"
$ws.Range("E14").Value = "  +7.74%  "
$ws.Range("D15").Value = "'0.635"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "33.906.06"
$ws.Range("E16").Value = "  +10.42%  "
$ws.Range("D17").Value = "'10.08"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "'4.24"
$ws.Range("E18").Value = "  +6.84%  "
$ws.Range("D19").Value = "'69.05"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("D20").Value = "'255.02"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("D21").Value = "0.0₃0739"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'10.36"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").Value = "'4.29"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").Value = "'158.89"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "'16.46"
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'7.03"
$ws.Range("E29").Value = "  +5.36%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +9.39%  "
$ws.Range("D32").Value = "'0.0507"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").Value = "'3.50"
$ws.Range("E34").Value = "  +7.00%  "
$ws.Range("D35").Value = "1.546.06"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").Value = "'1.80"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("D39").Value = "'83.64"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'0.615"
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").Value = "'2.83"
$ws.Range("E41").Value = "  +4.25%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "'0.899"
$ws.Range("E43").Value = "  +6.88%  "
$ws.Range("E44").Value = "  +5.48%  "
$ws.Range("D45").Value = "'0.0523"
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("D47").Value = "1.963.85"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("D50").Value = "'52.09"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("E51").Value = "  +6.83%  "

# The quote-prefix trick above marks the cell with a "quotePrefix"
# style; reset back to Normal so the cell style matches the original
# (unstyled) cells exactly.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
